$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.34"
$ws.Range("E2").Value = "'-1.40%"
$ws.Range("D3").Value = "'35.66"
$ws.Range("E3").Value = "'-0.73%"
$ws.Range("D4").Value = "'5.038"
$ws.Range("E4").Value = "'-0.13%"
$ws.Range("D5").Value = "'0.07971"
$ws.Range("E5").Value = "'-1.93%"
$ws.Range("D6").Value = "'1.848"
$ws.Range("E6").Value = "'-5.43%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.124"
$ws.Range("E7").Value = "'-0.18%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.760"
$ws.Range("E8").Value = "'-0.43%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9207"
$ws.Range("E9").Value = "'-1.26%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1268"
$ws.Range("E10").Value = "'-4.79%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1877"
$ws.Range("E11").Value = "'-2.07%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08936"
$ws.Range("E12").Value = "'-3.31%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03418"
$ws.Range("E13").Value = "'-2.51%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09856"
$ws.Range("E14").Value = "'-0.12%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001406"
$ws.Range("E15").Value = "'-0.59%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006287"
$ws.Range("E16").Value = "'8.51%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.865"
$ws.Range("E17").Value = "'7.29%"
$ws.Range("D18").Value = "'3.297"
$ws.Range("E18").Value = "'12.07%"
$ws.Range("D19").Value = "'0.3406"
$ws.Range("E19").Value = "'-0.71%"
$ws.Range("D20").Value = "'0.1341"
$ws.Range("E20").Value = "'0.67%"
$ws.Range("D21").Value = "'4.796"
$ws.Range("E21").Value = "'-7.57%"
$ws.Range("D22").Value = "'0.2346"
$ws.Range("E22").Value = "'-9.57%"
$ws.Range("D23").Value = "'0.04344"
$ws.Range("E23").Value = "'-0.81%"
$ws.Range("D24").Value = "'0.001237"
$ws.Range("E24").Value = "'1.39%"
$ws.Range("D25").Value = "'0.004837"
$ws.Range("E25").Value = "'1.39%"
$ws.Range("E27").Value = "'-21.11%"
$ws.Range("E39").Value = "'-3.89%"
$ws.Range("D40").Value = "'0.05095"
$ws.Range("E40").Value = "'0.63%"
$ws.Range("D41").Value = "'0.007581"
$ws.Range("E41").Value = "'-0.36%"
$ws.Range("E42").Value = "'-9.18%"
$ws.Range("D43").Value = "'0.1343"
$ws.Range("E43").Value = "'-2.57%"
$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'0.85%"
$ws.Range("E45").Value = "'-12.45%"
$ws.Range("D46").Value = "'0.00006203"
$ws.Range("E46").Value = "'-2.79%"
$ws.Range("E47").Value = "'0.51%"
$ws.Range("D48").Value = "'63.68"
$ws.Range("E48").Value = "'0.18%"
$ws.Range("D49").Value = "'0.001254"
$ws.Range("E49").Value = "'5.51%"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.51%"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.51%"
